$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1726
$ws.Range("F3").Value = 10068
$ws.Range("F8").Value = 1589
$ws.Range("F9").Value = 163
$ws.Range("F10").Value = 368
$ws.Range("F12").Value = 196
$ws.Range("F15").Value = 1168
$ws.Range("F19").Value = 83
$ws.Range("F21").Value = 13
$ws.Range("F23").Value = 96
$ws.Range("F24").Value = 1025
$ws.Range("F25").Value = 686
$ws.Range("F29").Value = 220
$ws.Range("F31").Value = 357
$ws.Range("F32").Value = 215
$ws.Range("F34").Value = 521
$ws.Range("F35").Value = 568
$ws.Range("F36").Value = 723
$ws.Range("F39").Value = 806
$ws.Range("F40").Value = 372
$ws.Range("F41").Value = 326

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F12").Value = 24
$ws.Range("F17").Value = 1076
$ws.Range("F19").Value = 568
$ws.Range("F20").Value = 1098
$ws.Range("F21").Value = 319
$ws.Range("F23").Value = 70
$ws.Range("F33").Value = 155
$ws.Range("G37").Value = 144

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value = 189
$ws.Range("F6").Value = 2501
$ws.Range("F7").Value = 4034
$ws.Range("F8").Value = 53
$ws.Range("F9").Value = 1
$ws.Range("F10").Value = 276
$ws.Range("F11").Value = 180

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 1726
$ws.Range("F4").Value = 10068
$ws.Range("F5").Value = 189
$ws.Range("F7").Value = 4034
$ws.Range("F8").Value = 53
$ws.Range("F9").Value = 276
$ws.Range("F10").Value = 276
$ws.Range("F12").Value = 1589
$ws.Range("F13").Value = 163
$ws.Range("F14").Value = 368
$ws.Range("F15").Value = 196
$ws.Range("F18").Value = 1168
$ws.Range("F21").Value = 24
$ws.Range("F24").Value = 83
$ws.Range("F25").Value = 1076
$ws.Range("F28").Value = 1098
$ws.Range("F29").Value = 319
$ws.Range("F30").Value = 1025
$ws.Range("F31").Value = 686
$ws.Range("F32").Value = 70
$ws.Range("F35").Value = 357
$ws.Range("F38").Value = 521
$ws.Range("F39").Value = 568
$ws.Range("F41").Value = 723
$ws.Range("F43").Value = 806
$ws.Range("F44").Value = 372
$ws.Range("F46").Value = 326
